# Updated cryptos list data (prices / 1h volume %) as refreshed by the
# scheduled GitHub Actions job. Rows 12 and 13 also swap which coin
# (WrappedEther / Polkadot) occupies that rank, along with their own
# updated price & volume figures.
#
# Values are prefixed with a literal leading apostrophe so Excel stores
# them as text (matching the workbook's original inlineStr cells) instead
# of auto-converting number-looking strings like "214.64" into floating
# point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''25.982.55'
$ws.Range("E2").Value = '''  +0.58%  '
$ws.Range("D3").Value = '''1.641.53'
$ws.Range("E3").Value = '''  +0.46%  '
$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '''  +0.32%  '
$ws.Range("D5").Value = '''214.64'
$ws.Range("E5").Value = '''  +0.49%  '
$ws.Range("D6").Value = '''0.5090'
$ws.Range("E6").Value = '''  +1.49%  '
$ws.Range("D7").Value = '''1.002'
$ws.Range("E7").Value = '''  +0.17%  '
$ws.Range("D8").Value = '''0.2561'
$ws.Range("E8").Value = '''  +0.26%  '
$ws.Range("D9").Value = '''0.06372'
$ws.Range("E9").Value = '''  +0.32%  '
$ws.Range("D10").Value = '''19.48'
$ws.Range("E10").Value = '''  +0.48%  '
$ws.Range("D11").Value = '''0.07773'
$ws.Range("E11").Value = '''  +0.21%  '
$ws.Range("B12").Value = '''Polkadot'
$ws.Range("C12").Value = '''https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = '''4.285'
$ws.Range("E12").Value = '''  +1.21%  '
$ws.Range("B13").Value = '''WrappedEther'
$ws.Range("C13").Value = '''https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '''1.649.78'
$ws.Range("E13").Value = '''  +0.22%  '
$ws.Range("D14").Value = '''0.5446'
$ws.Range("E14").Value = '''  +0.98%  '
$ws.Range("D15").Value = '''0.0₅7746'
$ws.Range("E15").Value = '''  -1.38%  '
$ws.Range("D16").Value = '''64.29'
$ws.Range("E16").Value = '''  +0.13%  '
$ws.Range("D17").Value = '''25.967.12'
$ws.Range("E17").Value = '''  +0.37%  '
$ws.Range("D18").Value = '''1.002'
$ws.Range("E18").Value = '''  -0.05%  '
$ws.Range("D19").Value = '''196.42'
$ws.Range("E19").Value = '''  +0.68%  '
$ws.Range("D20").Value = '''4.430'
$ws.Range("E20").Value = '''  +1.75%  '
$ws.Range("D21").Value = '''9.928'
$ws.Range("E21").Value = '''  +0.59%  '
$ws.Range("D22").Value = '''6.044'
$ws.Range("E22").Value = '''  +1.66%  '
$ws.Range("D23").Value = '''1.004'
$ws.Range("E23").Value = '''  +0.19%  '
$ws.Range("D24").Value = '''1.885'
$ws.Range("E24").Value = '''  +0.09%  '
$ws.Range("D25").Value = '''141.12'
$ws.Range("E25").Value = '''  +1.34%  '
$ws.Range("D26").Value = '''0.1197'
$ws.Range("E26").Value = '''  +5.69%  '
$ws.Range("D27").Value = '''6.849'
$ws.Range("E27").Value = '''  +0.74%  '
$ws.Range("D28").Value = '''15.61'
$ws.Range("E28").Value = '''  +0.14%  '
$ws.Range("D29").Value = '''1.235'
$ws.Range("E29").Value = '''  +0.19%  '
$ws.Range("D30").Value = '''0.04879'
$ws.Range("E30").Value = '''  +0.94%  '
$ws.Range("D31").Value = '''3.252'
$ws.Range("E31").Value = '''  +0.49%  '
$ws.Range("D32").Value = '''3.174'
$ws.Range("E32").Value = '''  +0.54%  '
$ws.Range("D33").Value = '''1.527'
$ws.Range("E33").Value = '''  +0.19%  '
$ws.Range("D34").Value = '''2.365'
$ws.Range("E34").Value = '''  +0.46%  '
$ws.Range("D35").Value = '''0.8938'
$ws.Range("E35").Value = '''  +1.21%  '
$ws.Range("D36").Value = '''1.146.76'
$ws.Range("E36").Value = '''  +2.24%  '
$ws.Range("D37").Value = '''2.580'
$ws.Range("E37").Value = '''  -0.48%  '
$ws.Range("D38").Value = '''0.5435'
$ws.Range("E38").Value = '''  -0.93%  '
$ws.Range("D39").Value = '''0.01557'
$ws.Range("E39").Value = '''  +0.24%  '
$ws.Range("D40").Value = '''1.002'
$ws.Range("E40").Value = '''  +0.17%  '
$ws.Range("D41").Value = '''2.522'
$ws.Range("E41").Value = '''  -1.81%  '
$ws.Range("D42").Value = '''0.0₈127'
$ws.Range("E42").Value = '''  +4.31%  '
$ws.Range("D43").Value = '''0.8127'
$ws.Range("E43").Value = '''  +0.37%  '
$ws.Range("D44").Value = '''99.01'
$ws.Range("E44").Value = '''  +0.03%  '
$ws.Range("D45").Value = '''5.439'
$ws.Range("E45").Value = '''  -3.68%  '
$ws.Range("D46").Value = '''1.777.24'
$ws.Range("E46").Value = '''  +0.23%  '
$ws.Range("D47").Value = '''0.4527'
$ws.Range("E47").Value = '''  +0.62%  '
$ws.Range("D48").Value = '''54.93'
$ws.Range("E48").Value = '''  +0.13%  '
$ws.Range("D49").Value = '''0.9972'
$ws.Range("E49").Value = '''  -0.89%  '
$ws.Range("D50").Value = '''0.05054'
$ws.Range("E50").Value = '''  +0.46%  '
$ws.Range("E51").Value = '''  -0.52%  '
